$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '243.49'

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '23.73'

# Row 4
$ws.Range("B4").Value = 'HuobiToken'
$ws.Range("C4").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '5.305'
$ws.Range("E4").Value = '3HuobiTokenHT'

# Row 5
$ws.Range("B5").Value = 'Cronos'
$ws.Range("C5").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.05794'
$ws.Range("E5").Value = '4CronosCRO'

# Row 6
$ws.Range("B6").Value = 'KuCoinToken'
$ws.Range("C6").Value = 'https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs'
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '6.493'
$ws.Range("E6").Value = '5KuCoinTokenKCS'

# Row 7
$ws.Range("B7").Value = 'GateToken'
$ws.Range("C7").Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '3.337'
$ws.Range("E7").Value = '6GateTokenGT'

# Row 8
$ws.Range("B8").Value = 'MXToken'
$ws.Range("C8").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.8090'
$ws.Range("E8").Value = '7MXTokenMX'

# Row 9
$ws.Range("B9").Value = 'FTXToken'
$ws.Range("C9").Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.8730'
$ws.Range("E9").Value = '8FTXTokenFTT'

# Row 10
$ws.Range("B10").Value = 'WazirX'
$ws.Range("C10").Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.1389'
$ws.Range("E10").Value = '9WazirXWRX'

# Row 11
$ws.Range("B11").Value = 'MandalaExchangeToken'
$ws.Range("C11").Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07278'
$ws.Range("E11").Value = '10MandalaExchangeTokenMDX'

# Row 12
$ws.Range("B12").Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range("C12").Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.03074'
$ws.Range("E12").Value = '11LiechtensteinCryptoassetsExchangeLCX'

# Row 13
$ws.Range("B13").Value = 'BitrueCoin'
$ws.Range("C13").Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.03057'
$ws.Range("E13").Value = '12BitrueCoinBTR'

# Row 14
$ws.Range("B14").Value = 'BitMartToken'
$ws.Range("C14").Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.09317'
$ws.Range("E14").Value = '13BitMartTokenBMX'

# Row 15
$ws.Range("B15").Value = 'MCDex'
$ws.Range("C15").Value = 'https://coinranking.com/coin/3nMM61qeg+mcdex-mcb'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '3.875'
$ws.Range("E15").Value = '14MCDexMCB'

# Row 16
$ws.Range("B16").Value = 'BitForexToken'
$ws.Range("C16").Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.001534'
$ws.Range("E16").Value = '15BitForexTokenBF'

# Row 17
$ws.Range("B17").Value = 'CoinExToken'
$ws.Range("C17").Value = 'https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.04720'
$ws.Range("E17").Value = '16CoinExTokenCET'

# Row 18
$ws.Range("B18").Value = 'One'
$ws.Range("C18").Value = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.0006053'
$ws.Range("E18").Value = '17OneONEWorstin24h'

# Row 19
$ws.Range("B19").Value = 'TigerCash'
$ws.Range("C19").Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.005970'
$ws.Range("E19").Value = '18TigerCashTCH'

# Row 20
$ws.Range("B20").Value = 'BitKan'
$ws.Range("C20").Value = 'https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.001270'
$ws.Range("E20").Value = '19BitKanKAN'

# Row 21
$ws.Range("B21").Value = 'HotbitToken'
$ws.Range("C21").Value = 'https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.004592'
$ws.Range("E21").Value = '20HotbitTokenHTB'

# Row 22
$ws.Range("B22").Value = 'NitroEx'
$ws.Range("C22").Value = 'https://coinranking.com/coin/8oiZw6gwYhC+nitroex-ntx'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.00008706'
$ws.Range("E22").Value = '21NitroExNTX'

# Row 23
$ws.Range("B23").Value = 'LEO'
$ws.Range("C23").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '3.577'
$ws.Range("E23").Value = '22LEOLEO'

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.144'

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.3208'

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.1319'

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.0002345'

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.03783'

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.003216'

# Row 42
$ws.Range("B42").Value = 'BKEXToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.1053'
$ws.Range("E42").Value = '41BKEXTokenBKK'

# Row 43
$ws.Range("B43").Value = 'CEJI'
$ws.Range("C43").Value = 'https://coinranking.com/coin/SbKjCVJCh+ceji-ceji'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.002462'
$ws.Range("E43").Value = '42CEJICEJI'

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.006993'

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.00005485'

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.5502'
$ws.Range("E47").Value = '46CoinbaseStockTokenCOIN'

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.006737'
$ws.Range("E48").Value = '47BOLOBOLOBestin24h'

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.00002101'

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0002001'
